$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.002.90'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '2.947.56'
$ws.Range('E3').Value = '  +2.96%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '553.68'
$ws.Range('E5').Value = '  +1.32%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '133.72'
$ws.Range('E6').Value = '  +10.86%  '
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('E8').Value = '  +5.58%  '
$ws.Range('D9').Value = '2.941.84'
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('E10').Value = '  +3.78%  '
$ws.Range('E11').Value = '  +1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.447'
$ws.Range('E12').Value = '  +5.37%  '
$ws.Range('E13').Value = '  +5.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.88'
$ws.Range('E14').Value = '  +6.41%  '
$ws.Range('E15').Value = '  +3.51%  '
$ws.Range('D16').Value = '3.432.11'
$ws.Range('E16').Value = '  +3.12%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '6.90'
$ws.Range('E17').Value = '  +10.55%  '
$ws.Range('D18').Value = '2.942.27'
$ws.Range('E18').Value = '  +3.12%  '
$ws.Range('D19').Value = '57.930.65'
$ws.Range('E19').Value = '  +1.49%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '417.66'
$ws.Range('E20').Value = '  +2.46%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.34'
$ws.Range('E21').Value = '  +5.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.695'
$ws.Range('E22').Value = '  +8.40%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.43'
$ws.Range('E23').Value = '  +9.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.01'
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '78.92'
$ws.Range('E25').Value = '  +4.07%  '
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('E27').Value = '  +0.05%  '
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.02'
$ws.Range('E29').Value = '  +7.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.49'
$ws.Range('E30').Value = '  +6.83%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '25.53'
$ws.Range('E31').Value = '  +4.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.93'
$ws.Range('E32').Value = '  +0.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0972'
$ws.Range('E33').Value = '  +5.94%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.70'
$ws.Range('E34').Value = '  +7.66%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.945'
$ws.Range('E35').Value = '  +7.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.09'
$ws.Range('E36').Value = '  +5.96%  '
$ws.Range('D37').Value = '0.0₃0699'
$ws.Range('E37').Value = '  +15.22%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '48.36'
$ws.Range('E38').Value = '  +0.40%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.75'
$ws.Range('E39').Value = '  +6.94%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.67'
$ws.Range('E40').Value = '  +15.09%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '381.83'
$ws.Range('E41').Value = '  +9.54%  '
$ws.Range('E42').Value = '  +4.68%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0348'
$ws.Range('E43').Value = '  +3.39%  '
$ws.Range('D44').Value = '2.701.11'
$ws.Range('E44').Value = '  +4.89%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '124.13'
$ws.Range('E46').Value = '  +6.40%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.237'
$ws.Range('E47').Value = '  +5.50%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.97'
$ws.Range('E48').Value = '  +4.24%  '
$ws.Range('E49').Value = '  +3.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.93'
$ws.Range('E50').Value = '  +3.12%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.01'
$ws.Range('E51').Value = '  +5.21%  '
